$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.557.79"
$ws.Range("E2").Value = "  +3.77%  "
$ws.Range("D3").Value = "2.421.53"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'315.06"
$ws.Range("E5").Value = "  +4.04%  "
$ws.Range("D6").Value = "'101.10"
$ws.Range("E6").Value = "  +5.85%  "
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = "  +8.21%  "
$ws.Range("D10").Value = "'35.37"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("D11").Value = "'0.0800"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").Value = "'19.00"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").Value = "'0.122"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").Value = "2.800.78"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").Value = "2.415.52"
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("E17").Value = "  +5.14%  "
$ws.Range("D18").Value = "44.437.52"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").Value = "'12.43"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("D21").Value = "0.0₃0926"
$ws.Range("E21").Value = "  +4.67%  "
$ws.Range("D22").Value = "'68.70"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").Value = "'243.07"
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("E24").Value = "  +5.60%  "
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'25.23"
$ws.Range("E27").Value = "  +3.17%  "
$ws.Range("E28").Value = "  -3.90%  "
$ws.Range("D29").Value = "'9.55"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").Value = "'33.08"
$ws.Range("E30").Value = "  +3.42%  "
$ws.Range("D31").Value = "'48.46"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").Value = "'0.125"
$ws.Range("E32").Value = "  +20.34%  "
$ws.Range("D33").Value = "'19.37"
$ws.Range("E33").Value = "  +11.00%  "
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("E35").Value = "  +7.72%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  +2.79%  "
$ws.Range("D38").Value = "'4.49"
$ws.Range("E38").Value = "  +4.35%  "
$ws.Range("D39").Value = "'2.86"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").Value = "'119.62"
$ws.Range("E42").Value = "  -6.56%  "
$ws.Range("D43").Value = "'20.98"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("D45").Value = "1.940.43"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "'2.17"
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "  +8.93%  "
$ws.Range("D48").Value = "'9.42"
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("D49").Value = "'1.67"
$ws.Range("E49").Value = "  +11.38%  "
$ws.Range("D50").Value = "'54.44"
$ws.Range("E50").Value = "  +6.08%  "
$ws.Range("D51").Value = "'75.27"
$ws.Range("E51").Value = "  +5.34%  "
